$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column widths -----------------------------------------------------
# The sheet's columns were widened (content got shorter numbers in H/I,
# so Excel's original authoring tool recalculated the autofit widths).
# ColumnWidth values below are chosen so the saved OOXML "width" comes out
# as close as possible to the target widths from the authored workbook.
$colWidths = @(32.833333333333336, 32.5, 31.0, 37.833333333333336, 31.166666666666668, 29.666666666666668, 36.666666666666664, 31.5, 30.166666666666668)
for ($i = 0; $i -lt $colWidths.Length; $i++) {
    $ws.Columns.Item($i + 1).ColumnWidth = $colWidths[$i]
}

# --- Data table (successful run for 0.5s timestep) ----------------------
# Row 2
$ws.Range("A2").Value = 0.0070000000000000001
$ws.Range("B2").Value = 33.207226185180005
$ws.Range("C2").Value = 16.270073689068269
$ws.Range("D2").Value = 0.48995581860219967
$ws.Range("E2").Value = 30.125609173907485
$ws.Range("F2").Value = 15.106612683022341
$ws.Range("G2").Value = 0.50145418125209373
$ws.Range("H2").Value = 282
$ws.Range("I2").Value = 248.5

# Row 3
$ws.Range("A3").Value = 0.0073333333333333332
$ws.Range("B3").Value = 33.028521380253522
$ws.Range("C3").Value = 16.156674875946859
$ws.Range("D3").Value = 0.48917342347654447
$ws.Range("E3").Value = 30.30614439007276
$ws.Range("F3").Value = 15.267554599718283
$ws.Range("G3").Value = 0.50377753115699542
$ws.Range("H3").Value = 281.5
$ws.Range("I3").Value = 248.5

# Row 4
$ws.Range("A4").Value = 0.0076666666666666671
$ws.Range("B4").Value = 33.189505294638941
$ws.Range("C4").Value = 16.248281182820662
$ws.Range("D4").Value = 0.48956081262968465
$ws.Range("E4").Value = 30.439520634747339
$ws.Range("F4").Value = 15.370079113649274
$ws.Range("G4").Value = 0.50493827738220065
$ws.Range("H4").Value = 281.5
$ws.Range("I4").Value = 248

# Row 5
$ws.Range("A5").Value = 0.0080000000000000002
$ws.Range("B5").Value = 33.328107027831287
$ws.Range("C5").Value = 16.32260635156851
$ws.Range("D5").Value = 0.48975497882126934
$ws.Range("E5").Value = 30.182555240022822
$ws.Range("F5").Value = 15.145552571550679
$ws.Range("G5").Value = 0.50179822255298312
$ws.Range("H5").Value = 282
$ws.Range("I5").Value = 247.5

# Row 6 (A value also changes, from 0.0086666666666666663 to 0.0083333333333333332,
# because the sweep points shifted down by one row from row 7 onward)
$ws.Range("A6").Value = 0.0083333333333333332
$ws.Range("B6").Value = 33.051379859167362
$ws.Range("C6").Value = 16.172682622016822
$ws.Range("D6").Value = 0.4893194381272119
$ws.Range("E6").Value = 30.292616485853422
$ws.Range("F6").Value = 15.243230345571105
$ws.Range("G6").Value = 0.50319952892447095
$ws.Range("H6").Value = 281
$ws.Range("I6").Value = 247.5

# Row 7
$ws.Range("A7").Value = 0.0086666666666666663
$ws.Range("B7").Value = 33.144376478583418
$ws.Range("C7").Value = 16.206766987962084
$ws.Range("D7").Value = 0.48897486421065289
$ws.Range("E7").Value = 30.380946399435807
$ws.Range("F7").Value = 15.305862172752631
$ws.Range("G7").Value = 0.50379807039312219
$ws.Range("H7").Value = 281.5
$ws.Range("I7").Value = 247.5

# Row 8
$ws.Range("A8").Value = 0.0090000000000000011
$ws.Range("B8").Value = 33.222573189920894
$ws.Range("C8").Value = 16.280306914143384
$ws.Range("D8").Value = 0.49003750615808783
$ws.Range("E8").Value = 30.420725389849224
$ws.Range("F8").Value = 15.326754929914344
$ws.Range("G8").Value = 0.50382608348414237
$ws.Range("H8").Value = 281.5
$ws.Range("I8").Value = 247

# Row 9
$ws.Range("A9").Value = 0.0093333333333333341
$ws.Range("B9").Value = 33.290167005814595
$ws.Range("C9").Value = 16.280345067674681
$ws.Range("D9").Value = 0.48904365859237325
$ws.Range("E9").Value = 30.495546772448243
$ws.Range("F9").Value = 15.390815272280866
$ws.Range("G9").Value = 0.50469058276357837
$ws.Range("H9").Value = 281.5
$ws.Range("I9").Value = 247

# Row 10
$ws.Range("A10").Value = 0.0096666666666666672
$ws.Range("B10").Value = 33.3482677238365
$ws.Range("C10").Value = 16.335566351619534
$ws.Range("D10").Value = 0.48984752332257675
$ws.Range("E10").Value = 30.525144095978675
$ws.Range("F10").Value = 15.394546399453459
$ws.Range("G10").Value = 0.50432346366815373
$ws.Range("H10").Value = 281.5
$ws.Range("I10").Value = 247

# Row 11
$ws.Range("A11").Value = 0.01
$ws.Range("B11").Value = 33.39462924544231
$ws.Range("C11").Value = 16.317367183436009
$ws.Range("D11").Value = 0.48862249865112661
$ws.Range("E11").Value = 30.590676489596902
$ws.Range("F11").Value = 15.480958121910632
$ws.Range("G11").Value = 0.50606785787085506
$ws.Range("H11").Value = 281.5
$ws.Range("I11").Value = 247
